# Natmi following Dr Hou advice
# Adds a third cluster ("ECs") to the LR-pair table, expanding the
# sending/target cluster combinations from 2x2 (4 rows) to 3x3 (9 rows),
# and refreshes every numeric metric accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column order: A Sending cluster, B Ligand symbol, C Receptor symbol,
# D Target cluster, E..T numeric metrics.
$ligand = "Inhba"
$receptor = "Acvr2b"

$rows = @(
    @{ Send="FAPs"; Target="FAPs"; Vals=@(2,0.6666666666666666,4.123204333333334,12.369613,0.2909967288544799,0.2909967288544799,2,0.6666666666666666,1.539665666666667,4.618997,0.3572088291809875,0.3572088291809875,6.348356148684557,57.13520533816101,0.103946600809606,0.103946600809606) },
    @{ Send="FAPs"; Target="sCs";  Vals=@(2,0.6666666666666666,4.123204333333334,12.369613,0.2909967288544799,0.2909967288544799,3,1,1.452872333333333,4.358617,0.3370724153789985,0.3370724153789985,5.990489500580112,53.914405505221,0.09808697026236705,0.09808697026236705) },
    @{ Send="FAPs"; Target="ECs";  Vals=@(2,0.6666666666666666,4.123204333333334,12.369613,0.2909967288544799,0.2909967288544799,3,1,1.317729666666667,3.953189,0.3057187554400141,0.3057187554400141,5.43326867176189,48.899418045857,0.08896315778250684,0.08896315778250684) },
    @{ Send="sCs";  Target="FAPs"; Vals=@(3,1,8.433639666666666,25.300919,0.5952073574179045,0.5952073574179045,2,0.6666666666666666,1.539665666666667,4.618997,0.3572088291809875,0.3572088291809875,12.98498543980478,116.864868958243,0.2126133232631592,0.2126133232631592) },
    @{ Send="sCs";  Target="sCs";  Vals=@(3,1,8.433639666666666,25.300919,0.5952073574179045,0.5952073574179045,3,1,1.452872333333333,4.358617,0.3370724153789985,0.3370724153789985,12.25300174100255,110.277015669023,0.200627981616204,0.200627981616204) },
    @{ Send="sCs";  Target="ECs";  Vals=@(3,1,8.433639666666666,25.300919,0.5952073574179045,0.5952073574179045,3,1,1.317729666666667,3.953189,0.3057187554400141,0.3057187554400141,11.11325718674344,100.019314680691,0.1819660525385414,0.1819660525385414) },
    @{ Send="ECs";  Target="FAPs"; Vals=@(3,1,1.612402333333333,4.837207,0.1137959137276156,0.1137959137276156,2,0.6666666666666666,1.539665666666667,4.618997,0.3572088291809875,0.3572088291809875,2.482560513486556,22.343044621379,0.04064890510822223,0.04064890510822222) },
    @{ Send="ECs";  Target="sCs";  Vals=@(3,1,1.612402333333333,4.837207,0.1137959137276156,0.1137959137276156,3,1,1.452872333333333,4.358617,0.3370724153789985,0.3370724153789985,2.342614740302111,21.083532662719,0.03835746350042752,0.03835746350042752) },
    @{ Send="ECs";  Target="ECs";  Vals=@(3,1,1.612402333333333,4.837207,0.1137959137276156,0.1137959137276156,3,1,1.317729666666667,3.953189,0.3057187554400141,0.3057187554400141,2.124710389235889,19.122393503123,0.03478954511896586,0.03478954511896585) }
)

$r = 2
foreach ($row in $rows) {
    $ws.Cells.Item($r, 1).Value = $row.Send
    $ws.Cells.Item($r, 2).Value = $ligand
    $ws.Cells.Item($r, 3).Value = $receptor
    $ws.Cells.Item($r, 4).Value = $row.Target

    $c = 5
    foreach ($v in $row.Vals) {
        $ws.Cells.Item($r, $c).Value = $v
        $c = $c + 1
    }
    $r = $r + 1
}
